$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-13 all share the same open/close/high/low/shares_outstanding values
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 4).Value = 15.44999980926514   # D: open_price
    $ws.Cells.Item($r, 5).Value = 1.350000023841858   # E: close_price
    $ws.Cells.Item($r, 6).Value = 15.5                # F: high_price
    $ws.Cells.Item($r, 7).Value = 1.350000023841858   # G: low_price
    $ws.Cells.Item($r, 8).Value = 0                   # H: shares_outstanding
}

# Row 14 has its own distinct values
$ws.Cells.Item(14, 4).Value = 0.5953999757766724
$ws.Cells.Item(14, 5).Value = 0.5723999738693237
$ws.Cells.Item(14, 6).Value = 0.7773000001907349
$ws.Cells.Item(14, 7).Value = 0.506600022315979
$ws.Cells.Item(14, 8).Value = 0

# Column I (fixed_ticker) for rows 2-14 all become "ZPTA" to match column A
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 9).Value = "ZPTA"
}
